$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new date + new "kotlet szwajcar" combo meal, price 33 -> 32
$ws.Range("A2").Value = "14.05.2025"
$ws.Range("B2").Value = "Kotlet „szwajcar”, ziemniaki, surówka + zupa jarzynowa lub fasolowa po bretońsku"
$ws.Range("C2").Value = 32

# Row 3: new date + new "potrawka z kurczaka" combo meal, price unchanged (30)
$ws.Range("A3").Value = "14.05.2025"
$ws.Range("B3").Value = "Potrawka z kurczaka w sosie porowym, ziemniaki, surówka + zupa jarzynowa lub fasolowa po bretońsku"

# Row 4: new date + new "kotlet szwajcar" single meal, price 30 -> 29
$ws.Range("A4").Value = "14.05.2025"
$ws.Range("B4").Value = "Kotlet ”szwajcar”, ziemniaki, surówka "
$ws.Range("C4").Value = 29

# Row 5: new date + new "potrawka z kurczaka" single meal, price unchanged (27)
$ws.Range("A5").Value = "14.05.2025"
$ws.Range("B5").Value = "Potrawka z kurczaka w sosie porowym, ziemniaki, surówka"

# Row 6: new date + new soup, price unchanged (9.5)
$ws.Range("A6").Value = "14.05.2025"
$ws.Range("B6").Value = "Zupa jarzynowa lub fasolowa po bretońsku"

# Rows 7-10: date updated only, text/prices unchanged
$ws.Range("A7").Value = "14.05.2025"
$ws.Range("A8").Value = "14.05.2025"
$ws.Range("A9").Value = "14.05.2025"
$ws.Range("A10").Value = "14.05.2025"

# Update the active cell selection to B10
[void]$ws.Range("B10").Select()
